# Added new Scripts in RCC Module.
# Adds a new test case row (RCC104) to the "Test Cases" sheet, mirroring
# the formatting of the existing data rows (row 10 is a plain, unmerged
# data row with the same style as most other rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 21

# Copy the formatting (styles/borders) of an existing plain data row so the
# new row matches the rest of the table, then paste only the formats.
$ws.Range("A10:E10").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new test case data.
$ws.Cells.Item($newRow, 1).Value = "RCC104"
$ws.Cells.Item($newRow, 2).Value = "ABCD4"
$ws.Cells.Item($newRow, 3).Value = "Verify that user is able to add an article to the group from search results  page."
$ws.Cells.Item($newRow, 4).Value = "Y"

# Match the author's final selection/cursor position.
$ws.Range("C21").Select()
